# Apply "traded, fixed issues with the repeater" edit:
#  - Fill in PriceChange/UpDown (X10/Y10) for the most recent existing row (row 10)
#  - Append a brand-new trading row (row 11) with a fresh date/verdict/etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 10 which previously had no PriceChange/UpDown values yet
$ws.Cells.Item(10, 24).Value = -0.010002000000000066   # X10 PriceChange
$ws.Cells.Item(10, 25).Value = "Down"                  # Y10 UpDown

# Add new row 11 for the next trading day
$ws.Cells.Item(11, 1).Value = 42654.882118055553   # A11 Date
$ws.Cells.Item(11, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(11, 2).Value = 20             # B11 ScoreFinal
$ws.Cells.Item(11, 3).Value = "Strong Buy"   # C11 Verdict
$ws.Cells.Item(11, 4).Value = 0              # D11 totalSentiment
$ws.Cells.Item(11, 5).Value = 0              # E11 wordCount
$ws.Cells.Item(11, 6).Value = 0              # F11 sentenceCount
$ws.Cells.Item(11, 7).Value = 0              # G11 posWordPercentage
$ws.Cells.Item(11, 8).Value = 0              # H11 negWordPercentage
$ws.Cells.Item(11, 9).Value = 0              # I11 posPhrasePercentage
$ws.Cells.Item(11, 10).Value = 0             # J11 negPhrasePercentage
$ws.Cells.Item(11, 11).Value = 0             # K11 ElapsedMs
$ws.Cells.Item(11, 12).Value = 0             # L11 posWordCount
$ws.Cells.Item(11, 13).Value = 0             # M11 negWordCount
$ws.Cells.Item(11, 14).Value = 0             # N11 positivePhraseCount
$ws.Cells.Item(11, 15).Value = 0             # O11 negativePhraseCount
$ws.Cells.Item(11, 16).Value = "Random"      # P11 Method
$ws.Cells.Item(11, 17).Value = 17.089518681678967   # Q11 RSI
$ws.Cells.Item(11, 18).Value = -24.44        # R11 PEG
$ws.Cells.Item(11, 19).Value = -0.1101       # S11 200Moving%
$ws.Cells.Item(11, 19).NumberFormat = "0.00%"
$ws.Cells.Item(11, 20).Value = -0.0419       # T11 50Moving%
$ws.Cells.Item(11, 20).NumberFormat = "0.00%"
$ws.Cells.Item(11, 21).Value = 6.47          # U11 PriceBook
$ws.Cells.Item(11, 22).Value = 1.88          # V11 Dividend
$ws.Cells.Item(11, 23).Value = -2            # W11 Bollinger
